# Quarterly indexing bug-fix: each date in column A currently marks the
# 1st day of the quarter-start month (Jan/Apr/Jul/Oct), but should mark the
# 15th day of the following month instead. Shift every date accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldSerial = $cell.Value2
    $oldDate = [DateTime]::FromOADate($oldSerial)
    $shifted = $oldDate.AddMonths(1)
    $newDate = Get-Date -Year $shifted.Year -Month $shifted.Month -Day 15 -Hour 0 -Minute 0 -Second 0
    $cell.Value2 = $newDate.ToOADate()
}
